$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.406.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.605.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.67%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.520'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.86%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.61'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.834.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.583.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.423.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.533'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("E20").Value = '  +3.78%  '
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  +2.19%  '
$ws.Range("E24").Value = '  +2.00%  '
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.60%  '
$ws.Range("E27").Value = '  +5.03%  '
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("E34").Value = '  +4.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.411.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("E38").Value = '  +4.69%  '
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  +2.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.536'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0491'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.49%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.796'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '52.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +22.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.745.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.77%  '
$ws.Range("E51").Value = '  -2.76%  '
